# "minor fix on cross layer"
# Adds a "blue" product-colour column + a "product" header to each of the
# three pairwise-comparison blocks, lowers the AHP importance weights in the
# block titles, and adds two new cross-tabulation (pairwise comparison)
# matrices per block under columns H:M that break down how each colour /
# product beat the others, with the AUS-NZ/NZ weighting notes off to the
# side.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Phase 1 - introduce brand-new shared strings in the exact order needed
# so the saved sharedStrings.xml table lines up with the target workbook:
#   6 blue, 7 product, 8 Australia(importance 1), 9 NZ(importance 0.5),
#   10 USA(importance 0.4), 11 AUS-NZ, 12 aus, 13 nz
# ---------------------------------------------------------------------
$ws.Range("E2").Value = "blue"
$ws.Range("F2").Value = "product"

$ws.Range("A1").Value = "Australia(importance 1)"
$ws.Range("A7").Value = "NZ(importance 0.5)"
$ws.Range("A13").Value = "USA(importance 0.4)"

$ws.Range("H1").Value = "AUS-NZ"
$ws.Range("M1").Value = "aus"
$ws.Range("H6").Value = "nz"

# ---------------------------------------------------------------------
# Phase 2 - fill in the rest of the "blue" column for the three existing
# colour tables (header already written above for the first block).
# ---------------------------------------------------------------------
$ws.Range("E8").Value = "blue"
$ws.Range("E14").Value = "blue"

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1
$ws.Range("E9").Value = 1
$ws.Range("E10").Value = 1
$ws.Range("E11").Value = 1
$ws.Range("E15").Value = 1
$ws.Range("E16").Value = 1
$ws.Range("E17").Value = 1

# ---------------------------------------------------------------------
# Phase 3 - new AUS-NZ pairwise-comparison matrix for the colours
# (silver / black / grey / blue), rows 1-5, columns H:M.
# ---------------------------------------------------------------------
$ws.Range("I1").Value = "silver"
$ws.Range("J1").Value = "black"
$ws.Range("K1").Value = "grey"
$ws.Range("L1").Value = "blue"

$ws.Range("H2").Value = "silver"
$ws.Range("I2").Value = 3
$ws.Range("J2").Value = 2
$ws.Range("K2").Value = 0
$ws.Range("L2").Value = 2

$ws.Range("H3").Value = "black"
$ws.Range("I3").Value = 2
$ws.Range("J3").Value = 1
$ws.Range("K3").Value = 1
$ws.Range("L3").Value = 4

$ws.Range("H4").Value = "grey"
$ws.Range("I4").Value = 1
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 2

$ws.Range("H5").Value = "blue"
$ws.Range("I5").Value = 2
$ws.Range("J5").Value = 1
$ws.Range("K5").Value = 0
$ws.Range("L5").Value = 1

# ---------------------------------------------------------------------
# Phase 4 - second AUS-NZ pairwise-comparison matrix, this time for the
# products (ipad / iphone / imac), rows 8-11, columns H:K.
# ---------------------------------------------------------------------
$ws.Range("I8").Value = "ipad"
$ws.Range("J8").Value = "iphone"
$ws.Range("K8").Value = "imac"

$ws.Range("H9").Value = "ipad"
$ws.Range("I9").Value = 1
$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 2

$ws.Range("H10").Value = "iphone"
$ws.Range("I10").Value = 2
$ws.Range("J10").Value = 3
$ws.Range("K10").Value = 4

$ws.Range("H11").Value = "imac"
$ws.Range("I11").Value = 1
$ws.Range("J11").Value = 1
$ws.Range("K11").Value = 5

# ---------------------------------------------------------------------
# Phase 5 - grey-out (explicit "no fill") the redundant lower-triangle +
# diagonal cells of both pairwise-comparison matrices, matching the
# formatting applied when the tables were built.
# ---------------------------------------------------------------------
$greyedOut = @("I2", "I3", "J3", "I4", "J4", "K4", "I5", "J5", "K5", "I10", "I11", "J11")
foreach ($ref in $greyedOut) {
    $cell = $ws.Range($ref)
    $cell.Interior.Color = 16777215
    $cell.Interior.Pattern = -4142
}

# ---------------------------------------------------------------------
# Phase 6 - match the saved selection / active cell.
# ---------------------------------------------------------------------
$ws.Range("H6").Select()
